# Fill in the "units" (F) and "population type" (G) columns for the
# per-parameter rows on the "Data inputs" sheet (rows 5-27) so they match
# the pattern already used by the rows below them (28-49), and update the
# sheet's saved scroll/selection state (frozen-pane scroll back to the top,
# selection moved to I6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data inputs")

$values = @{
    5  = @("years", "allpops")
    6  = @("years", "allpops")
    7  = @("years", "allpops")
    8  = @("years", "allpops")
    9  = @("years", "average")
    10 = @("years", "total")
    11 = @("years", "total")
    12 = @("years", "allpops")
    13 = @("years", "total")
    14 = @("years", "females")
    15 = @("years", "average")
    16 = @("years", "total")
    17 = @("years", "total")
    18 = @("years", "total")
    19 = @("years", "total")
    20 = @("years", "total")
    21 = @("years", "total")
    22 = @("years", "total")
    23 = @("years", "average")
    24 = @("years", "average")
    25 = @("years", "average")
    26 = @("years", "average")
    27 = @("years", "average")
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 6).Value = $pair[0]
    $ws.Cells.Item($row, 7).Value = $pair[1]
}

# Restore the view: scroll the frozen pane back to the top and select I6.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$ws.Range("I6").Select()
